$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 86
$ws.Range("K2").Value = 134
$ws.Range("F3").Value = 129
$ws.Range("G3").Value = 137
$ws.Range("H3").Value = 150
$ws.Range("I3").Value = 184
$ws.Range("J3").Value = 218
$ws.Range("K3").Value = 210
$ws.Range("D4").Value = 6
$ws.Range("C6").Value = 458
$ws.Range("D6").Value = 389
$ws.Range("F6").Value = 496
$ws.Range("I6").Value = 481
$ws.Range("J6").Value = 399
$ws.Range("C7").Value = 609
$ws.Range("D7").Value = 610
$ws.Range("F7").Value = 719
$ws.Range("G7").Value = 647
$ws.Range("H7").Value = 694
$ws.Range("I7").Value = 805
$ws.Range("J7").Value = 753
$ws.Range("K7").Value = 849

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("G3").Value = 6
$ws.Range("I6").Value = 31
$ws.Range("J6").Value = 27
$ws.Range("G7").Value = 46
$ws.Range("I7").Value = 47
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 4
$ws.Range("J5").Value = 12

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F2").Value = 8
$ws.Range("F7").Value = 55

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("H8").Value = 61
$ws.Range("I27").Value = 12
$ws.Range("F28").Value = 55
$ws.Range("K30").Value = 15
$ws.Range("G32").Value = 46
$ws.Range("I32").Value = 47
$ws.Range("J32").Value = 42
$ws.Range("K36").Value = 65
$ws.Range("J42").Value = 12
$ws.Range("D45").Value = 5
$ws.Range("F47").Value = 14
$ws.Range("I47").Value = 23
$ws.Range("I48").Value = 3
$ws.Range("C51").Value = 2
$ws.Range("D53").Value = 69
$ws.Range("J54").Value = 9
$ws.Range("K70").Value = 23
$ws.Range("G75").Value = 3
$ws.Range("F87").Value = 3
$ws.Range("J88").Value = 12
$ws.Range("C98").Value = 609
$ws.Range("D98").Value = 610
$ws.Range("F98").Value = 719
$ws.Range("G98").Value = 647
$ws.Range("H98").Value = 694
$ws.Range("I98").Value = 805
$ws.Range("J98").Value = 753
$ws.Range("K98").Value = 849

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("D6").Value = 41
$ws.Range("D7").Value = 69

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 3

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 10
$ws.Range("H5").Value = 12

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I4").Value = 7
$ws.Range("I5").Value = 12

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 2

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I3").Value = 1
$ws.Range("I6").Value = 3

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 3
$ws.Range("K6").Value = 15

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J4").Value = 7
$ws.Range("J5").Value = 9

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("D4").Value = 1
$ws.Range("D6").Value = 5

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("F3").Value = 4
$ws.Range("I5").Value = 18
$ws.Range("F6").Value = 14
$ws.Range("I6").Value = 23

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 4
$ws.Range("K6").Value = 23

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("G3").Value = 1
$ws.Range("G5").Value = 3

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("H3").Value = 8
$ws.Range("H7").Value = 61
